# Updates the crypto price-ticker data (cols D "Price" / E "Volume(1h)")
# with a fresh snapshot, and re-ranks rows 28-32 (PEPE moves up to rank 26,
# pushing WrappedeETH / Binance-PegBSC-USD / InternetComputer(DFINITY) /
# Fetch.AI down one slot each) to mirror the upstream coinranking.com feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold plain text in the source data (e.g. "1.00", "  +0.81%  "),
# not numbers. Excel's Range.Value setter auto-coerces numeric-looking
# strings ("1.00", "5.35", ...) into real numbers, which would corrupt the
# cell type versus the original inlineStr cells. Temporarily force the
# whole data range to Text format so every assignment below is stored
# verbatim as a string, then restore the default (General) formatting
# once all values are committed.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"


# Row 2
$ws.Range("D2").Value = "62.328.91"
$ws.Range("E2").Value = "  +1.15%  "

# Row 3
$ws.Range("D3").Value = "2.422.38"
$ws.Range("E3").Value = "  +1.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "562.52"
$ws.Range("E5").Value = "  +1.78%  "

# Row 6
$ws.Range("D6").Value = "143.83"
$ws.Range("E6").Value = "  +1.96%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").Value = "  +1.23%  "

# Row 9
$ws.Range("D9").Value = "2.419.81"
$ws.Range("E9").Value = "  +1.46%  "

# Row 10
$ws.Range("E10").Value = "  +1.34%  "

# Row 11
$ws.Range("E11").Value = "  -2.06%  "

# Row 12
$ws.Range("D12").Value = "5.35"

# Row 13
$ws.Range("E13").Value = "  +0.31%  "

# Row 14
$ws.Range("E14").Value = "  +0.67%  "

# Row 15
$ws.Range("E15").Value = "  +2.38%  "

# Row 16
$ws.Range("D16").Value = "2.860.03"
$ws.Range("E16").Value = "  +1.54%  "

# Row 17
$ws.Range("D17").Value = "62.116.55"
$ws.Range("E17").Value = "  +0.87%  "

# Row 18
$ws.Range("D18").Value = "2.422.97"
$ws.Range("E18").Value = "  +1.57%  "

# Row 19
$ws.Range("D19").Value = "11.31"
$ws.Range("E19").Value = "  +2.91%  "

# Row 20
$ws.Range("E20").Value = "  +0.88%  "

# Row 21
$ws.Range("D21").Value = "324.40"
$ws.Range("E21").Value = "  +0.47%  "

# Row 22
$ws.Range("E22").Value = "  +1.88%  "

# Row 23
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("D24").Value = "65.83"
$ws.Range("E24").Value = "  +2.40%  "

# Row 25
$ws.Range("E25").Value = "  -3.75%  "

# Row 26
$ws.Range("E26").Value = "  +1.57%  "

# Row 27
$ws.Range("D27").Value = "578.71"
$ws.Range("E27").Value = "  +8.41%  "

# Row 28
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0956"
$ws.Range("E28").Value = "  +4.64%  "

# Row 29
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.542.43"
$ws.Range("E29").Value = "  +1.57%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.47%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "8.28"
$ws.Range("E31").Value = "  +0.21%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.46"
$ws.Range("E32").Value = "  +3.69%  "

# Row 33
$ws.Range("E33").Value = "  +1.11%  "

# Row 34
$ws.Range("E34").Value = "  +2.09%  "

# Row 35
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +2.07%  "

# Row 36
$ws.Range("D36").Value = "5.72"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("D38").Value = "4.79"
$ws.Range("E38").Value = "  +0.61%  "

# Row 39
$ws.Range("E39").Value = "  +1.03%  "

# Row 40
$ws.Range("D40").Value = "152.53"
$ws.Range("E40").Value = "  +4.22%  "

# Row 41
$ws.Range("E41").Value = "  +0.57%  "

# Row 42
$ws.Range("E42").Value = "  -4.33%  "

# Row 43
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.12%  "

# Row 44
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").Value = "  +6.20%  "

# Row 45
$ws.Range("D45").Value = "150.07"
$ws.Range("E45").Value = "  +0.47%  "

# Row 46
$ws.Range("D46").Value = "3.66"
$ws.Range("E46").Value = "  +1.64%  "

# Row 47
$ws.Range("E47").Value = "  +2.06%  "

# Row 48
$ws.Range("D48").Value = "20.27"
$ws.Range("E48").Value = "  +1.33%  "

# Row 49
$ws.Range("D49").Value = "0.597"
$ws.Range("E49").Value = "  +2.14%  "

# Row 50
$ws.Range("D50").Value = "0.0923"
$ws.Range("E50").Value = "  +1.81%  "

# Row 51
$ws.Range("E51").Value = "  +2.24%  "

# Values are committed as text; drop the temporary Text number-format so
# the cells end up styled exactly as they started (no explicit format).
$dataRange.ClearFormats()
